$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("I3").Value = 3.5
$ws.Range("J3").Value = 1.08
$ws.Range("K3").Value = 7.5
$ws.Range("Z3").Value = 7.5
$ws.Range("AF3").Value = 17

# Row 18
$ws.Range("G18").Value = 1.95
$ws.Range("H18").Value = 3.3
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 1.07
$ws.Range("K18").Value = 9
$ws.Range("P18").Value = 1.44
$ws.Range("Q18").Value = 2.63
$ws.Range("R18").Value = 1.8
$ws.Range("S18").Value = 1.95
$ws.Range("T18").Value = 7
$ws.Range("U18").Value = 9
$ws.Range("Z18").Value = 9
$ws.Range("AA18").Value = 6.5
$ws.Range("AD18").Value = 251
$ws.Range("AF18").Value = 19
$ws.Range("AG18").Value = 13

# Row 19
$ws.Range("G19").Value = 1.5
$ws.Range("H19").Value = 4.33
$ws.Range("I19").Value = 6.25
$ws.Range("N19").Value = 1.93
$ws.Range("O19").Value = 1.93
$ws.Range("P19").Value = 1.4
$ws.Range("Q19").Value = 2.75
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 1.75
$ws.Range("T19").Value = 6.5
$ws.Range("W19").Value = 10
$ws.Range("AA19").Value = 8
$ws.Range("AE19").Value = 13
$ws.Range("AF19").Value = 29
$ws.Range("AG19").Value = 19
$ws.Range("AH19").Value = 67
